$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Set new value in B16 - adds a new shared string "static color (nach farbauswahl)"
$ws.Range("B16").Value = "static color (nach farbauswahl)"

# Update the active selection to B16 (as reflected in the diff's sheetView selection)
$ws.Range("B16").Select()
